$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 16570
$ws.Range("J3").Value = 16570
$ws.Range("L3").Value = 16570
$ws.Range("N3").Value = -16798
$ws.Range("H52").Value = 166669300
$ws.Range("I52").Value = 2933
$ws.Range("J52").Value = 333335680
$ws.Range("K52").Value = 8799
$ws.Range("L52").Value = 1000007040
$ws.Range("M52").Value = -8639
$ws.Range("N52").Value = -1000007360
$ws.Range("H99").Value = 285
$ws.Range("I99").Value = 285
$ws.Range("K99").Value = 855
$ws.Range("M99").Value = 643
$ws.Range("H102").Value = 16570
$ws.Range("J102").Value = 16570
$ws.Range("L102").Value = 16570
$ws.Range("N102").Value = -23060
$ws.Range("H113").Value = 5497521
$ws.Range("I113").Value = 20411650
$ws.Range("J113").Value = 2842.1052
$ws.Range("K113").Value = 20411650
$ws.Range("L113").Value = 2842.1052
$ws.Range("M113").Value = -20408396
$ws.Range("N113").Value = -9350.1052
$ws.Range("H138").Value = 4137.793
$ws.Range("I138").Value = 1585.0588
$ws.Range("J138").Value = 7754.1665
$ws.Range("K138").Value = 4755.1764
$ws.Range("L138").Value = 23262.4995
$ws.Range("M138").Value = 384.8235999999997
$ws.Range("N138").Value = -33542.49950000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 64281.9
$ws.Range("J23").Value = 62850.875
$ws.Range("L23").Value = 62850.875
$ws.Range("N23").Value = -63368.875
$ws.Range("H45").Value = 53508.21
$ws.Range("I45").Value = 143681.58
$ws.Range("J45").Value = 907.0833
$ws.Range("K45").Value = 143681.58
$ws.Range("L45").Value = 907.0833
$ws.Range("M45").Value = -143304.58
$ws.Range("N45").Value = -1661.0833
$ws.Range("H74").Value = 1479.0212
$ws.Range("I74").Value = 1077.1666
$ws.Range("J74").Value = 2188.1765
$ws.Range("K74").Value = 1077.1666
$ws.Range("L74").Value = 2188.1765
$ws.Range("M74").Value = -203.1666
$ws.Range("N74").Value = -3936.1765
$ws.Range("H77").Value = 1479.0212
$ws.Range("I77").Value = 1077.1666
$ws.Range("J77").Value = 2188.1765
$ws.Range("K77").Value = 5385.833000000001
$ws.Range("L77").Value = 10940.8825
$ws.Range("M77").Value = -1017.833000000001
$ws.Range("N77").Value = -19676.8825
$ws.Range("H102").Value = 2084.2727
$ws.Range("I102").Value = 1961
$ws.Range("K102").Value = 1961
$ws.Range("M102").Value = -339
$ws.Range("H122").Value = 1521.1562
$ws.Range("I122").Value = 1570.2413
$ws.Range("K122").Value = 4710.7239
$ws.Range("M122").Value = -2260.7239
$ws.Range("H132").Value = 1202565.1
$ws.Range("I132").Value = 1863.8055
$ws.Range("K132").Value = 5591.416499999999
$ws.Range("M132").Value = -3061.416499999999
$ws.Range("H134").Value = 46000
$ws.Range("J134").Value = 46000
$ws.Range("L134").Value = 46000
$ws.Range("N134").Value = -56140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 22728980
$ws.Range("I107").Value = 38462656
$ws.Range("J107").Value = 2559.889
$ws.Range("K107").Value = 38462656
$ws.Range("L107").Value = 2559.889
$ws.Range("M107").Value = -38460736
$ws.Range("N107").Value = -6399.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 47620380
$ws.Range("I58").Value = 71429620
$ws.Range("J58").Value = 1911.4286
$ws.Range("K58").Value = 71429620
$ws.Range("L58").Value = 1911.4286
$ws.Range("M58").Value = -71429417
$ws.Range("N58").Value = -2317.4286
$ws.Range("H62").Value = 5449.1
$ws.Range("J62").Value = 6061.5
$ws.Range("L62").Value = 6061.5
$ws.Range("N62").Value = -7309.5
$ws.Range("H65").Value = 5449.1
$ws.Range("J65").Value = 6061.5
$ws.Range("L65").Value = 30307.5
$ws.Range("N65").Value = -36547.5
$ws.Range("H99").Value = 62508810
$ws.Range("I99").Value = 500050000
$ws.Range("J99").Value = 2925.7144
$ws.Range("K99").Value = 500050000
$ws.Range("L99").Value = 2925.7144
$ws.Range("M99").Value = -500048502
$ws.Range("N99").Value = -5921.7144
$ws.Range("H126").Value = 62508810
$ws.Range("I126").Value = 500050000
$ws.Range("J126").Value = 2925.7144
$ws.Range("K126").Value = 1500150000
$ws.Range("L126").Value = 8777.143199999999
$ws.Range("M126").Value = -1500147530
$ws.Range("N126").Value = -13717.1432
$ws.Range("H132").Value = 15153508
$ws.Range("I132").Value = 1656.8334
$ws.Range("J132").Value = 83336840
$ws.Range("K132").Value = 4970.5002
$ws.Range("L132").Value = 250010520
$ws.Range("M132").Value = -2440.5002
$ws.Range("N132").Value = -250015580
$ws.Range("H136").Value = 47620380
$ws.Range("I136").Value = 71429620
$ws.Range("J136").Value = 1911.4286
$ws.Range("K136").Value = 214288860
$ws.Range("L136").Value = 5734.2858
$ws.Range("M136").Value = -214286310
$ws.Range("N136").Value = -10834.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20.82353
$ws.Range("I2").Value = 17.285715
$ws.Range("J2").Value = 37.333332
$ws.Range("K2").Value = 103.71429
$ws.Range("L2").Value = 223.999992
$ws.Range("M2").Value = 9.285709999999995
$ws.Range("N2").Value = -449.999992
$ws.Range("H9").Value = 137525580
$ws.Range("I9").Value = 33333700
$ws.Range("J9").Value = 200040720
$ws.Range("K9").Value = 100001100
$ws.Range("L9").Value = 600122160
$ws.Range("M9").Value = -100000876
$ws.Range("N9").Value = -600122608
$ws.Range("H34").Value = 527.6
$ws.Range("I34").Value = 347
$ws.Range("K34").Value = 1041
$ws.Range("M34").Value = -957
$ws.Range("H39").Value = 7425
$ws.Range("J39").Value = 7425
$ws.Range("L39").Value = 22275
$ws.Range("N39").Value = -22863
$ws.Range("H102").Value = 4893.75
$ws.Range("I102").Value = 4787.5
$ws.Range("K102").Value = 14362.5
$ws.Range("M102").Value = -11928.5
$ws.Range("H105").Value = 5200
$ws.Range("J105").Value = 5200
$ws.Range("L105").Value = 15600
$ws.Range("N105").Value = -20842
$ws.Range("H131").Value = 821.85
$ws.Range("J131").Value = 842.7826
$ws.Range("L131").Value = 2528.3478
$ws.Range("N131").Value = -12608.3478

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2150.4688
$ws.Range("I102").Value = 1712.1154
$ws.Range("K102").Value = 1712.1154
$ws.Range("M102").Value = -90.11539999999991
$ws.Range("H122").Value = 28433590
$ws.Range("I122").Value = 41555164
$ws.Range("J122").Value = 3516.3333
$ws.Range("K122").Value = 124665492
$ws.Range("L122").Value = 10548.9999
$ws.Range("M122").Value = -124663042
$ws.Range("N122").Value = -15448.9999
$ws.Range("H126").Value = 4207.25
$ws.Range("I126").Value = 3882.4
$ws.Range("J126").Value = 4354.909
$ws.Range("K126").Value = 11647.2
$ws.Range("L126").Value = 13064.727
$ws.Range("M126").Value = -9177.200000000001
$ws.Range("N126").Value = -18004.727
$ws.Range("H132").Value = 5976.4053
$ws.Range("I132").Value = 2874.423
$ws.Range("J132").Value = 13308.363
$ws.Range("K132").Value = 8623.269
$ws.Range("L132").Value = 39925.089
$ws.Range("M132").Value = -6093.269
$ws.Range("N132").Value = -44985.089

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2425
$ws.Range("I7").Value = 1900
$ws.Range("K7").Value = 1900
$ws.Range("M7").Value = -1788
$ws.Range("H40").Value = 50003880
$ws.Range("I40").Value = 2400
$ws.Range("J40").Value = 62504250
$ws.Range("K40").Value = 2400
$ws.Range("L40").Value = 62504250
$ws.Range("M40").Value = -2264
$ws.Range("N40").Value = -62504522
$ws.Range("H122").Value = 7277.6113
$ws.Range("I122").Value = 9314.77
$ws.Range("J122").Value = 1981
$ws.Range("K122").Value = 27944.31
$ws.Range("L122").Value = 5943
$ws.Range("M122").Value = -25494.31
$ws.Range("N122").Value = -10843
$ws.Range("H126").Value = 2425
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -3230
$ws.Range("H132").Value = 16672662
$ws.Range("I132").Value = 29413962
$ws.Range("J132").Value = 10961.846
$ws.Range("K132").Value = 88241886
$ws.Range("L132").Value = 32885.538
$ws.Range("M132").Value = -88239356
$ws.Range("N132").Value = -37945.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 371.75
$ws.Range("I107").Value = 460.66666
$ws.Range("K107").Value = 1381.99998
$ws.Range("M107").Value = 538.0000199999999
$ws.Range("H126").Value = 2679.0667
$ws.Range("I126").Value = 2440.0833
$ws.Range("J126").Value = 3635
$ws.Range("K126").Value = 7320.249899999999
$ws.Range("L126").Value = 10905
$ws.Range("M126").Value = -4850.249899999999
$ws.Range("N126").Value = -15845
$ws.Range("H132").Value = 84807.42999999999
$ws.Range("I132").Value = 118589.78
$ws.Range("J132").Value = 23999.2
$ws.Range("K132").Value = 355769.34
$ws.Range("L132").Value = 71997.60000000001
$ws.Range("M132").Value = -353239.34
$ws.Range("N132").Value = -77057.60000000001
